$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 301, shifting existing rows 301:329 down to 302:330
$ws.Rows("301:301").Insert()

# Populate the new row 301 with the new weekly data point for Membrillo (Vega Modelo de Temuco)
$ws.Range("A301").Value = 10
$ws.Range("B301").Value = "Vega Modelo de Temuco"
$ws.Range("C301").Value = "La Araucanía"
$ws.Range("D301").Value = 45166
$ws.Range("E301").Value = 9
$ws.Range("F301").Value = "Fruta"
$ws.Range("G301").Value = 100104
$ws.Range("H301").Value = "Frutos de pepita"
$ws.Range("I301").Value = 100104003
$ws.Range("J301").Value = "Membrillo"
$ws.Range("K301").Value = "Champion"
$ws.Range("L301").Value = "Primera"
$ws.Range("M301").Value = 155
$ws.Range("N301").Value = 16000
$ws.Range("O301").Value = 16000
$ws.Range("P301").Value = 16000
$ws.Range("Q301").Value = "$/bandeja 18 kilos granel"
$ws.Range("R301").Value = "Región de O'Higgins"
$ws.Range("S301").Value = 889
$ws.Range("T301").Value = 18
